$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 204.36363
$ws.Range("I53").Value = 186.4
$ws.Range("J53").Value = 219.33333
$ws.Range("K53").Value = 186.4
$ws.Range("L53").Value = 219.33333
$ws.Range("M53").Value = 450.6
$ws.Range("N53").Value = -1493.33333

$ws.Range("H70").Value = 5723.077
$ws.Range("I70").Value = 7737.5
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 23212.5
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -22942.5
$ws.Range("N70").Value = -8040

$ws.Range("H73").Value = 5723.077
$ws.Range("I73").Value = 7737.5
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 23212.5
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -22276.5
$ws.Range("N73").Value = -9372

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = 0

$ws.Range("H132").Value = 2195.7778
$ws.Range("I132").Value = 1186.5428
$ws.Range("J132").Value = 5728.1
$ws.Range("K132").Value = 3559.6284
$ws.Range("L132").Value = 17184.3
$ws.Range("M132").Value = -1029.6284
$ws.Range("N132").Value = -22244.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1972.3636
$ws.Range("I2").Value = 2010.7778
$ws.Range("J2").Value = 1799.5
$ws.Range("K2").Value = 2010.7778
$ws.Range("L2").Value = 1799.5
$ws.Range("M2").Value = -1897.7778
$ws.Range("N2").Value = -2025.5

$ws.Range("H61").Value = 502719.03
$ws.Range("I61").Value = 501671.56
$ws.Range("J61").Value = 503766.5
$ws.Range("K61").Value = 501671.56
$ws.Range("L61").Value = 503766.5
$ws.Range("M61").Value = -501459.56
$ws.Range("N61").Value = -504190.5

$ws.Range("H102").Value = 4178.148
$ws.Range("I102").Value = 1300.5264
$ws.Range("J102").Value = 11012.5
$ws.Range("K102").Value = 1300.5264
$ws.Range("L102").Value = 11012.5
$ws.Range("M102").Value = 321.4736
$ws.Range("N102").Value = -14256.5

$ws.Range("H116").Value = 1972.3636
$ws.Range("I116").Value = 2010.7778
$ws.Range("J116").Value = 1799.5
$ws.Range("K116").Value = 2010.7778
$ws.Range("L116").Value = 1799.5
$ws.Range("M116").Value = 283.2221999999999
$ws.Range("N116").Value = -6387.5

$ws.Range("H136").Value = 502719.03
$ws.Range("I136").Value = 501671.56
$ws.Range("J136").Value = 503766.5
$ws.Range("K136").Value = 1505014.68
$ws.Range("L136").Value = 1511299.5
$ws.Range("M136").Value = -1502464.68
$ws.Range("N136").Value = -1516399.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1972.3636
$ws.Range("I3").Value = 2010.7778
$ws.Range("J3").Value = 1799.5
$ws.Range("K3").Value = 2010.7778
$ws.Range("L3").Value = 1799.5
$ws.Range("M3").Value = -1896.7778
$ws.Range("N3").Value = -2027.5

$ws.Range("H105").Value = 7694031
$ws.Range("I105").Value = 1485.7142
$ws.Range("J105").Value = 16668667
$ws.Range("K105").Value = 1485.7142
$ws.Range("L105").Value = 16668667
$ws.Range("M105").Value = 261.2858000000001
$ws.Range("N105").Value = -16672161

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2677.525
$ws.Range("I31").Value = 1461.6072
$ws.Range("J31").Value = 5514.6665
$ws.Range("K31").Value = 1461.6072
$ws.Range("L31").Value = 5514.6665
$ws.Range("M31").Value = -1166.6072
$ws.Range("N31").Value = -6104.6665

$ws.Range("H34").Value = 2677.525
$ws.Range("I34").Value = 1461.6072
$ws.Range("J34").Value = 5514.6665
$ws.Range("K34").Value = 1461.6072
$ws.Range("L34").Value = 5514.6665
$ws.Range("M34").Value = -1259.6072
$ws.Range("N34").Value = -5918.6665

$ws.Range("H99").Value = 54680.844
$ws.Range("I99").Value = 85249.336
$ws.Range("J99").Value = 2277.7144
$ws.Range("K99").Value = 85249.336
$ws.Range("L99").Value = 2277.7144
$ws.Range("M99").Value = -83751.336
$ws.Range("N99").Value = -5273.7144

$ws.Range("H122").Value = 2805.8462
$ws.Range("I122").Value = 3052.3635
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 9157.0905
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -6707.0905
$ws.Range("N122").Value = -9250

$ws.Range("H126").Value = 54680.844
$ws.Range("I126").Value = 85249.336
$ws.Range("J126").Value = 2277.7144
$ws.Range("K126").Value = 255748.008
$ws.Range("L126").Value = 6833.1432
$ws.Range("M126").Value = -253278.008
$ws.Range("N126").Value = -11773.1432

$ws.Range("H132").Value = 2095.2683
$ws.Range("I132").Value = 1263.909
$ws.Range("J132").Value = 3057.8948
$ws.Range("K132").Value = 3791.727
$ws.Range("L132").Value = 9173.6844
$ws.Range("M132").Value = -1261.727
$ws.Range("N132").Value = -14233.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.033333
$ws.Range("I12").Value = 28.375
$ws.Range("J12").Value = 48.363636
$ws.Range("K12").Value = 85.125
$ws.Range("L12").Value = 145.090908
$ws.Range("M12").Value = 87.875
$ws.Range("N12").Value = -491.090908

$ws.Range("H98").Value = 3235.8823
$ws.Range("I98").Value = 226.5
$ws.Range("J98").Value = 4161.846
$ws.Range("K98").Value = 679.5
$ws.Range("L98").Value = 12485.538
$ws.Range("M98").Value = 818.5
$ws.Range("N98").Value = -15481.538

$ws.Range("H100").Value = 7916.0713
$ws.Range("I100").Value = 4025
$ws.Range("J100").Value = 8215.385
$ws.Range("K100").Value = 12075
$ws.Range("L100").Value = 24646.155
$ws.Range("M100").Value = -11264
$ws.Range("N100").Value = -26268.155

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = 0

$ws.Range("H102").Value = 6722.769
$ws.Range("I102").Value = 4750
$ws.Range("J102").Value = 7081.4546
$ws.Range("K102").Value = 14250
$ws.Range("L102").Value = 21244.3638
$ws.Range("M102").Value = -11816
$ws.Range("N102").Value = -26112.3638

$ws.Range("H103").Value = 2211.1904
$ws.Range("I103").Value = 856.4286
$ws.Range("J103").Value = 2888.5715
$ws.Range("K103").Value = 2569.2858
$ws.Range("L103").Value = 8665.7145
$ws.Range("M103").Value = -1690.2858
$ws.Range("N103").Value = -10423.7145

$ws.Range("H104").Value = 16900.5
$ws.Range("I104").Value = 1403
$ws.Range("J104").Value = 20000
$ws.Range("K104").Value = 4209
$ws.Range("L104").Value = 60000
$ws.Range("M104").Value = -1588
$ws.Range("N104").Value = -65242

$ws.Range("H105").Value = 6558.4443
$ws.Range("I105").Value = 3026
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 9078
$ws.Range("L105").Value = 21000
$ws.Range("M105").Value = -6457
$ws.Range("N105").Value = -26242

$ws.Range("H106").Value = 4411.45
$ws.Range("I106").Value = 3000
$ws.Range("J106").Value = 4485.737
$ws.Range("K106").Value = 9000
$ws.Range("L106").Value = 13457.211
$ws.Range("M106").Value = -8054
$ws.Range("N106").Value = -15349.211

$ws.Range("H113").Value = 19231378
$ws.Range("I113").Value = 26316376
$ws.Range("J113").Value = 671.4286
$ws.Range("K113").Value = 78949128
$ws.Range("L113").Value = 2014.2858
$ws.Range("M113").Value = -78946958
$ws.Range("N113").Value = -6354.2858

$ws.Range("H136").Value = 2633.8333
$ws.Range("I136").Value = 1385.3077
$ws.Range("J136").Value = 5880
$ws.Range("K136").Value = 4155.9231
$ws.Range("L136").Value = 17640
$ws.Range("M136").Value = 944.0769
$ws.Range("N136").Value = -27840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2928.6155
$ws.Range("I126").Value = 2797.75
$ws.Range("J126").Value = 3138
$ws.Range("K126").Value = 8393.25
$ws.Range("L126").Value = 9414
$ws.Range("M126").Value = -5923.25
$ws.Range("N126").Value = -14354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 45458380
$ws.Range("I2").Value = 4250
$ws.Range("J2").Value = 50003796
$ws.Range("K2").Value = 4250
$ws.Range("L2").Value = 50003796
$ws.Range("M2").Value = -4138
$ws.Range("N2").Value = -50004020

$ws.Range("H22").Value = 865.46155
$ws.Range("I22").Value = 612.75
$ws.Range("J22").Value = 977.7778
$ws.Range("K22").Value = 612.75
$ws.Range("L22").Value = 977.7778
$ws.Range("M22").Value = -317.75
$ws.Range("N22").Value = -1567.7778

$ws.Range("H27").Value = 865.46155
$ws.Range("I27").Value = 612.75
$ws.Range("J27").Value = 977.7778
$ws.Range("K27").Value = 612.75
$ws.Range("L27").Value = 977.7778
$ws.Range("M27").Value = -505.75
$ws.Range("N27").Value = -1191.7778

$ws.Range("H100").Value = 62504996
$ws.Range("I100").Value = 6674.5
$ws.Range("J100").Value = 166668860
$ws.Range("K100").Value = 6674.5
$ws.Range("L100").Value = 166668860
$ws.Range("M100").Value = -6133.5
$ws.Range("N100").Value = -166669942

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5438.4614
$ws.Range("I96").Value = 2929.6
$ws.Range("J96").Value = 13801.333
$ws.Range("K96").Value = 2929.6
$ws.Range("L96").Value = 13801.333
$ws.Range("M96").Value = -1556.6
$ws.Range("N96").Value = -16547.333

$ws.Range("H100").Value = 4756.2666
$ws.Range("I100").Value = 1863.4286
$ws.Range("J100").Value = 7287.5
$ws.Range("K100").Value = 3726.8572
$ws.Range("L100").Value = 14575
$ws.Range("M100").Value = -3185.8572
$ws.Range("N100").Value = -15657
